$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'43.749.35"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.49%  "

# Row 3
$ws.Range("D3").Value = "'2.287.81"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.15%  "

# Row 4
$ws.Range("E4").Value = "  +0.43%  "

# Row 5
$ws.Range("D5").Value = "'110.41"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +15.69%  "

# Row 6
$ws.Range("D6").Value = "'267.90"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.00%  "

# Row 7
$ws.Range("D7").Value = "'0.625"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.74%  "

# Row 8
$ws.Range("E8").Value = "  +0.30%  "

# Row 9
$ws.Range("D9").Value = "'0.616"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.45%  "

# Row 10
$ws.Range("D10").Value = "'47.69"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.32%  "

# Row 11
$ws.Range("D11").Value = "'0.0947"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.60%  "

# Row 12
$ws.Range("D12").Value = "'9.17"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +15.99%  "

# Row 13
$ws.Range("E13").Value = "  +0.58%  "

# Row 14
$ws.Range("D14").Value = "'15.81"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.58%  "

# Row 15
$ws.Range("D15").Value = "'2.629.18"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.11%  "

# Row 16
$ws.Range("D16").Value = "'0.847"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.80%  "

# Row 17
$ws.Range("D17").Value = "'2.279.92"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.22%  "

# Row 18
$ws.Range("D18").Value = "'43.589.19"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.18%  "

# Row 19
$ws.Range("D19").Value = "'0.0000110"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.51%  "

# Row 20
$ws.Range("D20").Value = "'6.84"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +10.28%  "

# Row 21
$ws.Range("D21").Value = "'72.17"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.26%  "

# Row 22
$ws.Range("E22").Value = "  -4.38%  "

# Row 23
$ws.Range("D23").Value = "'232.32"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.08%  "

# Row 24
$ws.Range("D24").Value = "'9.77"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +7.24%  "

# Row 25
$ws.Range("D25").Value = "'2.77"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +8.36%  "

# Row 26
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").Value = "'11.76"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +5.12%  "

# Row 27
$ws.Range("B27").Value = "Dai"
$ws.Range("C27").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D27").Value = "'0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.04%  "

# Row 28
$ws.Range("B28").Value = "InjectiveProtocol"
$ws.Range("C28").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D28").Value = "'41.63"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.00%  "

# Row 29
$ws.Range("B29").Value = "WEMIXToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D29").Value = "'3.39"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.25%  "

# Row 30
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "'2.27"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.55%  "

# Row 31
$ws.Range("B31").Value = "Monero"
$ws.Range("C31").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D31").Value = "'175.47"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.27%  "

# Row 32
$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").Value = "'0.0928"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.96%  "

# Row 33
$ws.Range("D33").Value = "'21.55"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.97%  "

# Row 34
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").Value = "'5.63"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.23%  "

# Row 35
$ws.Range("B35").Value = "Stellar"
$ws.Range("C35").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D35").Value = "'0.127"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.39%  "

# Row 36
$ws.Range("B36").Value = "RenderToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D36").Value = "'4.70"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +7.77%  "

# Row 37
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").Value = "'0.0366"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.95%  "

# Row 38
$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").Value = "'0.107"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.57%  "

# Row 39
$ws.Range("B39").Value = "NEARProtocol"
$ws.Range("C39").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D39").Value = "'3.80"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +13.97%  "

# Row 40
$ws.Range("B40").Value = "Algorand"
$ws.Range("C40").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D40").Value = "'0.244"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.27%  "

# Row 41
$ws.Range("D41").Value = "'13.82"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +12.66%  "

# Row 42
$ws.Range("B42").Value = "LidoDAOToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D42").Value = "'2.40"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.20%  "

# Row 43
$ws.Range("B43").Value = "MultiversX"
$ws.Range("C43").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D43").Value = "'72.96"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +11.14%  "

# Row 44
$ws.Range("B44").Value = "THORChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D44").Value = "'6.22"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +20.71%  "

# Row 45
$ws.Range("B45").Value = "FirstDigitalUSD"
$ws.Range("C45").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D45").Value = "'1.00"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.19%  "

# Row 46
$ws.Range("B46").Value = "ARBITRUM"
$ws.Range("C46").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D46").Value = "'1.38"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.53%  "

# Row 47
$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").Value = "'8.82"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.80%  "

# Row 48
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").Value = "'102.52"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +7.03%  "

# Row 49
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "'0.0992"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.43%  "

# Row 50
$ws.Range("B50").Value = "TrustWalletToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D50").Value = "'1.23"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.29%  "

# Row 51
$ws.Range("B51").Value = "WOONetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$ws.Range("D51").Value = "'0.451"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.14%  "
